$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the new "Anwesend" (present) column K for rows 5-9, matching the
# green "Anwesend" styling already used on columns I/J (copy value+format
# from J5, which already carries the correct fill).
$ws.Range("J5").Copy($ws.Range("K5"))
$ws.Range("J5").Copy($ws.Range("K6"))
$ws.Range("J5").Copy($ws.Range("K7"))
$ws.Range("J5").Copy($ws.Range("K8"))
$ws.Range("J5").Copy($ws.Range("K9"))

$excel.CutCopyMode = $false

# Matches the final selection recorded in the saved workbook.
$ws.Range("K9").Select()
